$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): columns B..L ---
$headers = @(
  "Algorithm",
  "One Year Base mean",
  "One Year Base std",
  "Two Year Base mean",
  "Two Year Base std",
  "Three Year Base mean",
  "Three Year Base std",
  "Five Year Base mean",
  "Five Year Base std",
  "Ten Year Base mean",
  "Ten Year Base std"
)
for ($i = 0; $i -lt $headers.Count; $i++) {
  $ws.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

# Newly added header cells H1:L1 need the same "header" style (bold, thin
# border, centered) that B1:G1 already carry. Copy the format from the
# existing G1 header cell instead of rebuilding it property-by-property.
$ws.Range("G1").Copy()
foreach ($col in @("H","I","J","K","L")) {
  $ws.Range($col + "1").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- Data rows 2..8: col A index, col B algorithm name, col C..L metric values ---
$rows = @(
  @(0, "LR",    0.8569752827817343, 0.02040235051356787, 0.8407979407979408, 0.02396534764999987, 0.8304215363100589, 0.02730114188011467, 0.8165092054263565, 0.04698609276682866, 0.8277983920841064, 0.03821836381816607),
  @(1, "LDA",   0.8453330540427315, 0.02226636598904934, 0.8279279279279279, 0.03238048814100637, 0.8106482868844271, 0.02789325696747806, 0.7963117732558139, 0.05198510947394241, 0.8095444238301381, 0.04081626926861007),
  @(2, "KNN",   0.8925722664432341, 0.02233268711618139, 0.8868404118404118, 0.02060136985134997, 0.888332833882729,  0.02181979368469639, 0.8833999515503876, 0.02206614330032326, 0.8581735724592867, 0.02988061301590509),
  @(3, "DTREE", 0.7838542103058231, 0.04244395587088937, 0.7981016731016732, 0.0316331609890594,  0.7738287883328339, 0.03758951121426484, 0.779251453488372,  0.0527492860491276,  0.7487012987012986, 0.04775212846402369),
  @(4, "RTREE", 0.8919145370758275, 0.0219666514902566,  0.8847812097812098, 0.02307479082431547, 0.8827240035960443, 0.01954645038142135, 0.8818616763565892, 0.04288578487513406, 0.9026592455163884, 0.0310281674216055),
  @(5, "XTREE", 0.8724842899036448, 0.02650603053825204, 0.8611141753998897, 0.02841045785277064, 0.8558485665767657, 0.02654332386504935, 0.853875968992248,  0.04079864627162615, 0.8460214388785818, 0.04003138132180345),
  @(6, "SVM",   0.8718516966904064, 0.02155964354370622, 0.8678801250229821, 0.03236166552384019, 0.8608280891019879, 0.02691800836013376, 0.8554324127906977, 0.03482333217432464, 0.8561533704390847, 0.04650982689853566)
)

$r = 2
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  for ($c = 0; $c -lt 10; $c++) {
    $ws.Cells.Item($r, 3 + $c).Value = $row[2 + $c]
  }
  $r++
}

# --- Remove the old trailing row 9 (sheet now ends at row 8) ---
$ws.Range("A9:L9").Delete()
